# Trade #13 closed at 2026-02-17 20:52:11 - unknown UNKNOWN +0.000%
#
# This updates the workbook to reflect the closing of the open
# MarketMaking trade (row 42 on "All Trades" / row 9 on "MarketMaking"),
# and rolls the new totals up into the "Summary" and "Strategy Status"
# sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet - overall portfolio totals
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.41   # Current Capital
$summary.Range("B4").Value = 0.19      # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 41        # Total Trades
$summary.Range("B7").Value = 18        # Winning Trades
$summary.Range("B9").Value = 43.9      # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking strategy row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.41     # Capital
$status.Range("D5").Value = 8          # Trades
$status.Range("E5").Value = 0.08       # P&L $
$status.Range("F5").Value = 0.41       # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------
# All Trades sheet - trade #41, row 42
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G42").Value = 0.130822   # Exit Price
$allTrades.Range("H42").Value = "CLOSED"   # Status
$allTrades.Range("I42").Value = 336.0719   # P&L %
$allTrades.Range("J42").Value = 0.1        # P&L $
$allTrades.Range("K42").Value = 100.41     # Capital After
$allTrades.Range("L42").Value = "early_exit" # Exit Reason
$allTrades.Range("M42").Value = 2.86       # Duration (min)

# ---------------------------------------------------------------
# MarketMaking sheet - same trade #41, row 9
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G9").Value = 0.130822   # Exit Price
$mm.Range("H9").Value = "CLOSED"   # Status
$mm.Range("I9").Value = 336.0719   # P&L %
$mm.Range("J9").Value = 0.1        # P&L $
$mm.Range("K9").Value = 100.41     # Capital After
$mm.Range("P9").Value = "early_exit" # Exit Reason
$mm.Range("Q9").Value = 2.86       # Duration (min)
